# Applies the cryptos.xlsx data-refresh diff via Excel COM interop.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The cells below hold plain decimal-looking text (e.g. "1.00", "607.08")
# that Excel would otherwise auto-coerce into a real number, dropping the
# significant trailing zeros the source data relies on. Mark each as Text
# first so the assigned string is preserved verbatim, matching the
# original inline-string cell content.
foreach ($addr in @(
    "D5",
    "D6",
    "D8",
    "D9",
    "D12",
    "D13",
    "D14",
    "D17",
    "D19",
    "D20",
    "D21",
    "D22",
    "D23",
    "D24",
    "D25",
    "D26",
    "D27",
    "D31",
    "D32",
    "D33",
    "D34",
    "D35",
    "D36",
    "D37",
    "D38",
    "D39",
    "D40",
    "D41",
    "D42",
    "D43",
    "D45",
    "D46",
    "D47",
    "D48",
    "D49",
    "D50"
    )) {
    $ws.Range($addr).NumberFormat = "@"
}

# Coin / Link / Price / Volume(1h) updates
$ws.Range("D2").Value = "72.976.97"
$ws.Range("E2").Value = "  +4.50%  "
$ws.Range("D3").Value = "2.646.57"
$ws.Range("E3").Value = "  +5.02%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "607.08"
$ws.Range("E5").Value = "  +1.92%  "
$ws.Range("D6").Value = "179.45"
$ws.Range("E6").Value = "  +1.07%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "0.528"
$ws.Range("E8").Value = "  +1.86%  "
$ws.Range("D9").Value = "0.173"
$ws.Range("E9").Value = "  +10.55%  "
$ws.Range("D10").Value = "2.644.74"
$ws.Range("E10").Value = "  +4.91%  "
$ws.Range("E11").Value = "  +0.84%  "
$ws.Range("D12").Value = "0.354"
$ws.Range("E12").Value = "  +3.55%  "
$ws.Range("D13").Value = "5.08"
$ws.Range("E13").Value = "  +1.45%  "
$ws.Range("D14").Value = "0.0000190"
$ws.Range("E14").Value = "  +6.97%  "
$ws.Range("D15").Value = "3.138.21"
$ws.Range("E15").Value = "  +6.13%  "
$ws.Range("D16").Value = "72.918.78"
$ws.Range("D17").Value = "26.83"
$ws.Range("E17").Value = "  +3.19%  "
$ws.Range("D18").Value = "2.658.74"
$ws.Range("E18").Value = "  +7.31%  "
$ws.Range("D19").Value = "386.80"
$ws.Range("E19").Value = "  +6.30%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "8.05"
$ws.Range("E20").Value = "  +5.39%  "
$ws.Range("B21").Value = "Chainlink"
$ws.Range("C21").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D21").Value = "11.57"
$ws.Range("E21").Value = "  +4.60%  "
$ws.Range("D22").Value = "4.20"
$ws.Range("E22").Value = "  +3.67%  "
$ws.Range("D23").Value = "2.02"
$ws.Range("E23").Value = "  +20.50%  "
$ws.Range("D24").Value = "73.26"
$ws.Range("E24").Value = "  +3.81%  "
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("B26").Value = "NEARProtocol"
$ws.Range("C26").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D26").Value = "4.40"
$ws.Range("E26").Value = "  +4.06%  "
$ws.Range("D27").Value = "9.90"
$ws.Range("E27").Value = "  +8.91%  "
$ws.Range("D28").Value = "2.789.53"
$ws.Range("E28").Value = "  +5.51%  "
$ws.Range("E29").Value = "  +0.43%  "
$ws.Range("D30").Value = "0.0₃0968"
$ws.Range("E30").Value = "  +8.35%  "
$ws.Range("D31").Value = "534.36"
$ws.Range("E31").Value = "  +4.19%  "
$ws.Range("D32").Value = "8.06"
$ws.Range("E32").Value = "  +3.69%  "
$ws.Range("D33").Value = "1.34"
$ws.Range("E33").Value = "  +8.56%  "
$ws.Range("D34").Value = "1.84"
$ws.Range("E34").Value = "  +3.40%  "
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.17%  "
$ws.Range("D36").Value = "163.55"
$ws.Range("E36").Value = "  +1.01%  "
$ws.Range("D37").Value = "19.37"
$ws.Range("E37").Value = "  +2.97%  "
$ws.Range("D38").Value = "0.113"
$ws.Range("E38").Value = "  -5.22%  "
$ws.Range("D39").Value = "1.41"
$ws.Range("E39").Value = "  +7.82%  "
$ws.Range("D40").Value = "19.11"
$ws.Range("E40").Value = "  +2.18%  "
$ws.Range("D41").Value = "1.85"
$ws.Range("E41").Value = "  +6.83%  "
$ws.Range("D42").Value = "5.12"
$ws.Range("E42").Value = "  +6.42%  "
$ws.Range("D43").Value = "2.64"
$ws.Range("E43").Value = "  +13.16%  "
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("D45").Value = "0.335"
$ws.Range("E45").Value = "  +4.48%  "
$ws.Range("D46").Value = "39.72"
$ws.Range("E46").Value = "  +2.27%  "
$ws.Range("D47").Value = "152.03"
$ws.Range("E47").Value = "  +1.32%  "
$ws.Range("D48").Value = "3.70"
$ws.Range("E48").Value = "  +3.35%  "
$ws.Range("D49").Value = "0.545"
$ws.Range("E49").Value = "  +5.81%  "
$ws.Range("D50").Value = "1.71"
$ws.Range("E50").Value = "  +8.64%  "
$ws.Range("D51").Value = "0.0₆0266"
$ws.Range("E51").Value = "  +6.48%  "
